$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 36; this pushes the existing rows 36-46 down to 37-47
# and expands the used range to A1:T47.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly price entry.
$ws.Range("A36").Value = 5
$ws.Range("B36").Value = "Macroferia Regional de Talca"
$ws.Range("C36").Value = "Maule"
$ws.Range("D36").Value = 45089
$ws.Range("E36").Value = 7
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100107
$ws.Range("H36").Value = "Otros"
$ws.Range("I36").Value = 100107011
$ws.Range("J36").Value = "Tuna"
$ws.Range("K36").Value = "Sin especificar"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 30
$ws.Range("N36").Value = 22000
$ws.Range("O36").Value = 22000
$ws.Range("P36").Value = 22000
$ws.Range("Q36").Value = '$/caja 18 kilos'
$ws.Range("R36").Value = "Provincia de Melipilla"
$ws.Range("S36").Value = 1222
$ws.Range("T36").Value = 18
